# Add a "See source of core program..." sentence (with a hyperlink to the
# Matlab cal_settling.html page) to the empty paragraph that follows the
# "Bibliography" heading, mirroring the existing sentence that already
# appears later in the document (near "Program Structure").

$d = $word.ActiveDocument

# Locate the empty paragraph right after the "Bibliography" heading and
# right before "Acknowledgements".
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim().Length -eq 0) {
        $prevText = ""
        if ($i -gt 1) { $prevText = $d.Paragraphs.Item($i - 1).Range.Text.Trim() }
        $nextText = ""
        if ($i -lt $d.Paragraphs.Count) { $nextText = $d.Paragraphs.Item($i + 1).Range.Text.Trim() }
        if ($prevText -eq "Bibliography" -and $nextText -eq "Acknowledgements") {
            $target = $para
            break
        }
    }
}

# Insertion point: start of the (empty) paragraph, before its end-of-
# paragraph mark, so the new runs land inside the existing <w:p>.
$ins = $d.Range($target.Range.Start, $target.Range.Start)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:iCs/>
              </w:rPr>
              <w:t xml:space="preserve">See source of core program and further details at: </w:t>
            </w:r>
            <w:hyperlink r:id="hlinkCalSettlingBiblio" w:history="1">
              <w:r>
                <w:rPr>
                  <w:rStyle w:val="Hyperlink"/>
                  <w:iCs/>
                </w:rPr>
                <w:t>http://neumeier.perso.ch/matlab/cal_settling.html</w:t>
              </w:r>
            </w:hyperlink>
            <w:r>
              <w:rPr>
                <w:iCs/>
                <w:color w:val="E6A82F" w:themeColor="accent1"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="hlinkCalSettlingBiblio" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://neumeier.perso.ch/matlab/cal_settling.html" TargetMode="External"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$ins.InsertXML($xml)
